$wb = $excel.ActiveWorkbook

$oldName = "F-SW-SD-02"
$newName = "F-SW-SD-03"

$ws = $wb.Worksheets.Item($oldName)
$ws.Name = $newName

foreach ($n in $wb.Names) {
    if ($n.RefersTo -like "*'$oldName'*") {
        $n.RefersTo = $n.RefersTo -replace [regex]::Escape("'$oldName'"), "'$newName'"
    }
}
